$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/custom style previously applied to A1 (back to default formatting)
$ws.Range("A1").ClearFormats()

# Append the two new words to the list
$ws.Range("A5").Value = "hitler"
$ws.Range("A6").Value = "queen"

# Move the active selection to A7, matching where the user clicked next
$null = $ws.Range("A7").Select()
